# Add a new "BOUNDARY" attack block (columns AK:AR) for SEED 88,
# mirroring the structure of the existing attack blocks (REV, REV_BIM,
# FGSM_SURRO, FGSM) in row 1/2, and fill in the corresponding MAE / RMSE /
# SIM metric values for the LSTM, RNN and GRU model rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("AK","AL","AM","AN","AO","AP","AQ","AR")

# ---------------------------------------------------------------------
# Row 1: header label "BOUNDARY" spanning AK1:AR1 (mirrors AC1 "FGSM" etc)
# (merge first, then apply formatting/value so the merge doesn't end up
# splitting the shared border style across the merged cells)
# ---------------------------------------------------------------------
$ws.Range("AK1:AR1").Merge()
$ws.Range("AJ1").Copy()
$ws.Range("AK1:AR1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("AK1").Value = "BOUNDARY"

# ---------------------------------------------------------------------
# Row 2: epsilon labels 0.01 .. 0.20 (mirrors AC2:AJ2 etc)
# ---------------------------------------------------------------------
$ws.Range("AJ2").Copy()
$ws.Range("AK2:AR2").PasteSpecial(-4122)   # xlPasteFormats

$epsLabels = @("0.01", "0.02", "0.03", "0.04", "0.05", "0.07", "0.10", "0.20")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $epsLabels[$i]
}

# ---------------------------------------------------------------------
# New metric values for the BOUNDARY attack, columns AK:AR, rows 4-12
# ---------------------------------------------------------------------
$newData = @{
    4  = @(390.3992577044169, 391.8017908287048, 398.8038006528218, 401.2359029515584, 421.3227803293864, 437.3160461616516, 463.2117540168762, 651.9611461575826)
    5  = @(489.2113279360275, 491.5688788238218, 494.9450467217215, 497.7337667428206, 530.339194723864,  554.1363678360606, 567.5598559118322, 818.6349930379963)
    6  = @(0.9990755579861272, 0.9990673890700759, 0.9990557359144541, 0.9990432672346851, 0.9989105669777159, 0.9988262245787846, 0.9987713311790573, 0.9975022087184466)
    7  = @(537.1415200042725, 540.281365292867,  541.9926981226603, 555.7811170260112, 558.6827701505025, 571.9447099622091, 595.4838717969259, 775.259315624237)
    8  = @(649.671967341643,  653.7355689588217, 658.1318471867484, 678.5656354234328, 673.1440240242024, 694.3897178746944, 717.6490616161614, 927.8302652424078)
    9  = @(0.9990671685500947, 0.9990470537408496, 0.9990643617599748, 0.9989418088728429, 0.9989692877497833, 0.9988450356825828, 0.9986796817364582, 0.9972872214855639)
    10 = @(310.8977380053202, 315.2886210187276, 320.1887220573425, 323.6158444023133, 335.4404909070333, 382.5507218424479, 409.1945012728373, 676.986574529012)
    11 = @(428.7031259102286, 431.3301281492065, 435.3153895389307, 439.1204777305023, 448.764046995636,  509.2524243842903, 528.005104241012,  885.9279588019787)
    12 = @(0.9993372423833846, 0.9993306313821481, 0.9993179189799252, 0.9993042057978951, 0.9992698019827521, 0.9990401874922721, 0.9989535001256086, 0.9969671125511961)
}

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# Tiny recomputation deltas on pre-existing SIM cells (rows 6, 9, 12)
# ---------------------------------------------------------------------
$modifications = @(
    @{ Cell = "F6";  Value = 0.9990687356509985 }
    @{ Cell = "J6";  Value = 0.9990089146365944 }
    @{ Cell = "K6";  Value = 0.9989420615148382 }
    @{ Cell = "N6";  Value = 0.9990822397610548 }
    @{ Cell = "R6";  Value = 0.9990747589810025 }
    @{ Cell = "S6";  Value = 0.9990665717275093 }
    @{ Cell = "W6";  Value = 0.9981172722701451 }
    @{ Cell = "Y6";  Value = 0.9971664375423135 }
    @{ Cell = "AC6"; Value = 0.9987903087061284 }
    @{ Cell = "AD6"; Value = 0.9984271428906181 }
    @{ Cell = "AE6"; Value = 0.9979740729933866 }
    @{ Cell = "AI6"; Value = 0.9932168311616112 }
    @{ Cell = "E9";  Value = 0.9990202767977286 }
    @{ Cell = "I9";  Value = 0.9986559571695052 }
    @{ Cell = "M9";  Value = 0.9990202767977286 }
    @{ Cell = "P9";  Value = 0.9987956391672492 }
    @{ Cell = "U9";  Value = 0.9989322763309333 }
    @{ Cell = "V9";  Value = 0.9987562371960799 }
    @{ Cell = "AC9"; Value = 0.9988943855504203 }
    @{ Cell = "AD9"; Value = 0.9986661050700377 }
    @{ Cell = "AE9"; Value = 0.9983825694895385 }
    @{ Cell = "AG9"; Value = 0.9976122245473198 }
    @{ Cell = "C12"; Value = 0.9993477063791676 }
    @{ Cell = "D12"; Value = 0.999237468333524 }
    @{ Cell = "P12"; Value = 0.9992314012229793 }
    @{ Cell = "R12"; Value = 0.9990714537975353 }
    @{ Cell = "X12"; Value = 0.9984518659263646 }
    @{ Cell = "Y12"; Value = 0.9981222684962211 }
    @{ Cell = "AA12"; Value = 0.995599645285814 }
    @{ Cell = "AE12"; Value = 0.9985937376995941 }
)

foreach ($m in $modifications) {
    $ws.Range($m.Cell).Value = $m.Value
}

Write-Host "BOUNDARY attack block added for SEED 88"
